# feature/ADMINDASH-354: change template from anwari
#
# Source workbook: client/assets/template/templateImportOrders.xlsx
#
# Changes applied:
#   1. "Order list" sheet, cell S2 (PickupType column): REGULER -> LATER
#   2. "Order list" sheet, cell AA2 (PaymentType column): WALLET -> Wallet
#   3. "legend" sheet, cell C2 (Payment Type legend): WALLET -> Wallet
#   4. "legend" sheet, cell C3 (Payment Type legend): CASH -> Cash
#   5. "Order list" sheet active selection moves from R2 to U4
#
# Dropping the old literal strings "WALLET"/"CASH" from every cell that used
# them causes the shared-string table to lose those two now-unused entries
# and gain the newly-typed "Wallet"/"Cash" strings on save.

$wb = $excel.ActiveWorkbook

$orders = $wb.Worksheets.Item("Order list")
$legend = $wb.Worksheets.Item("legend")

# Data corrections on the "Order list" sheet (row 2 = sample order)
$orders.Range("S2").Value = "LATER"
$orders.Range("AA2").Value = "Wallet"

# Matching casing fix on the "legend" sheet lookup list used for data validation
$legend.Range("C2").Value = "Wallet"
$legend.Range("C3").Value = "Cash"

# Move the active selection on the "Order list" sheet to U4
$orders.Activate() | Out-Null
$orders.Range("U4").Select() | Out-Null
